$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
        @("ECs","Fn1","Itgb6","ECs",3,1,19.95578266666667,59.867348,0.0117373419656925,0.0117373419656925,2,0.6666666666666666,0.3530253333333334,1.059076,0.05211323950342944,0.05211323950342944,7.044896827827556,63.40407145044801,0.0006116709129917864,0.0006116709129917864),
        @("ECs","Fn1","Itgb6","FAPs",3,1,19.95578266666667,59.867348,0.0117373419656925,0.0117373419656925,3,1,2.310677666666666,6.932033,0.3410998794937062,0.3410998794937063,46.11138132872044,415.002431958484,0.004003605930074132,0.004003605930074133),
        @("ECs","Fn1","Itgb6","M2",3,1,19.95578266666667,59.867348,0.0117373419656925,0.0117373419656925,1,0.3333333333333333,0.07663533333333333,0.229906,0.01131282971314188,0.01131282971314188,1.529318056587556,13.763862509288,0.0001327825509427932,0.0001327825509427932),
        @("ECs","Fn1","Itgb6","sCs",3,1,19.95578266666667,59.867348,0.0117373419656925,0.0117373419656925,3,1,4.033858333333334,12.101575,0.5954740512897225,0.5954740512897225,80.49880020812223,724.4892018731,0.006989282571683787,0.006989282571683787),
        @("FAPs","Fn1","Itgb6","ECs",3,1,1637.343343333333,4912.03003,0.9630320723052701,0.9630320723052702,2,0.6666666666666666,0.3530253333333334,1.059076,0.05211323950342944,0.05211323950342944,578.0236795613645,5202.21311605228,0.05018672103352852,0.05018672103352852),
        @("FAPs","Fn1","Itgb6","FAPs",3,1,1637.343343333333,4912.03003,0.9630320723052701,0.9630320723052702,3,1,2.310677666666666,6.932033,0.3410998794937062,0.3410998794937063,3783.372696105665,34050.35426495099,0.3284901238119018,0.3284901238119019),
        @("FAPs","Fn1","Itgb6","M2",3,1,1637.343343333333,4912.03003,0.9630320723052701,0.9630320723052702,1,0.3333333333333333,0.07663533333333333,0.229906,0.01131282971314188,0.01131282971314188,125.4783528974644,1129.30517607718,0.01089461784228366,0.01089461784228366),
        @("FAPs","Fn1","Itgb6","sCs",3,1,1637.343343333333,4912.03003,0.9630320723052701,0.9630320723052702,3,1,4.033858333333334,12.101575,0.5954740512897225,0.5954740512897225,6604.811090033028,59443.29981029725,0.5734606096175562,0.5734606096175563),
        @("M2","Fn1","Itgb6","ECs",3,1,17.50081933333334,52.502458,0.01029341242216722,0.01029341242216722,2,0.6666666666666666,0.3530253333333334,1.059076,0.05211323950342944,0.05211323950342944,6.178232578756446,55.60409320880801,0.0005364230668639759,0.0005364230668639761),
        @("M2","Fn1","Itgb6","FAPs",3,1,17.50081933333334,52.502458,0.01029341242216722,0.01029341242216722,3,1,2.310677666666666,6.932033,0.3410998794937062,0.3410998794937063,40.43875238190156,363.948771437114,0.003511081736780257,0.003511081736780258),
        @("M2","Fn1","Itgb6","M2",3,1,17.50081933333334,52.502458,0.01029341242216722,0.01029341242216722,1,0.3333333333333333,0.07663533333333333,0.229906,0.01131282971314188,0.01131282971314188,1.341181123216445,12.070630108948,0.000116447621899117,0.000116447621899117),
        @("M2","Fn1","Itgb6","sCs",3,1,17.50081933333334,52.502458,0.01029341242216722,0.01029341242216722,3,1,4.033858333333334,12.101575,0.5954740512897225,0.5954740512897225,70.59582590792779,635.3624331713501,0.006129459996623869,0.00612945999662387),
        @("sCs","Fn1","Itgb6","ECs",3,1,25.39612333333333,76.18836999999999,0.01493717330687017,0.01493717330687017,2,0.6666666666666666,0.3530253333333334,1.059076,0.05211323950342944,0.05211323950342944,8.965474905124445,80.68927414612,0.0007784244900451581,0.0007784244900451582),
        @("sCs","Fn1","Itgb6","FAPs",3,1,25.39612333333333,76.18836999999999,0.01493717330687017,0.01493717330687017,3,1,2.310677666666666,6.932033,0.3410998794937062,0.3410998794937063,58.68225500624555,528.1402950562099,0.00509506801495002,0.005095068014950021),
        @("sCs","Fn1","Itgb6","M2",3,1,25.39612333333333,76.18836999999999,0.01493717330687017,0.01493717330687017,1,0.3333333333333333,0.07663533333333333,0.229906,0.01131282971314188,0.01131282971314188,1.946240377024444,17.51616339322,0.0001689816980163106,0.0001689816980163106),
        @("sCs","Fn1","Itgb6","sCs",3,1,25.39612333333333,76.18836999999999,0.01493717330687017,0.01493717330687017,3,1,4.033858333333334,12.101575,0.5954740512897225,0.5954740512897225,102.4443637425278,921.99927368275,0.008894699103858681,0.008894699103858681)
)

$r = 2
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    for ($c = 5; $c -le 20; $c++) {
        $ws.Cells.Item($r, $c).Value = $row[$c - 1]
    }
    $r++
}
